# Remove the <w:contextualSpacing w:val="0"/> element from every paragraph's
# properties (w:pPr) in the document body. The Word object model does not
# expose a ContextualSpacing property on ParagraphFormat, so we operate on
# the raw OOXML for the whole document body via WordOpenXML / InsertXML.

$d = $word.ActiveDocument
$r = $d.Content

$xml = $r.WordOpenXML
$newXml = $xml -replace '\s*<w:contextualSpacing[^/]*/>', ''

# Self-closing empty elements (e.g. <w:rPr/>) can get silently dropped by
# the InsertXML re-parser, which would lose paragraph mark run-properties
# that were not part of this edit. Spell them out as open/close pairs so
# they survive the round-trip untouched.
$newXml = $newXml -replace '<(w:\w+)\s*/>', '<$1></$1>'

if ($newXml -ne $xml) {
    $r.InsertXML($newXml)
}
